# Applies the "Week 16 logged + season sim from Week 17" update to the
# Browns Team Data workbook.
#
# 1) Appends the Week-16 per-play logs to the four running play logs on the
#    YDS sheet (rush/pass yards gained, for Offense and Defense).
# 2) Appends the simulated special-teams logs (Week 17+) to the six running
#    logs on the ST sheet (KO/PT distance, return attempts, return yards).
# 3) Updates the season-to-date summary totals on OFF, DEF, ST, TURNS and
#    PEN sheets to reflect the newly logged/simulated weeks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append Week 16 play-by-play yardage logs
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 1 1 1 27 11 17 3 8 3 12 5 8 13 8 15 13 3 3 12 8 30 4 10 4"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 7 40 12 8 8 11 24 14 1 16 5 3 10 10 12 3 11 10 5 7"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 0 4 8 2 13 8 3 2 3 3 27 6 6 6 3 10 7 1 -1 6 5 4"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 1 8 4 33 4 6 11 34 1 9 2 -2 10 9 12 1 10 1 17 10 1 9 7 4"

# ---------------------------------------------------------------------
# OFF sheet - update season totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 14
$offWs.Range("C2").Value = 473
$offWs.Range("D2").Value = 22
$offWs.Range("F2").Value = 133
$offWs.Range("G2").Value = 100
$offWs.Range("I2").Value = 15
$offWs.Range("J2").Value = 68
$offWs.Range("N2").Value = 37
$offWs.Range("O2").Value = 53
$offWs.Range("P2").Value = 24

$offWs.Range("B3").Value = 24
$offWs.Range("C3").Value = 289
$offWs.Range("D3").Value = 15
$offWs.Range("E3").Value = 64
$offWs.Range("F3").Value = 212
$offWs.Range("G3").Value = 68
$offWs.Range("H3").Value = 56
$offWs.Range("I3").Value = 100
$offWs.Range("J3").Value = 94
$offWs.Range("L3").Value = 557
$offWs.Range("M3").Value = 351
$offWs.Range("Q3").Value = 1050

# ---------------------------------------------------------------------
# DEF sheet - update season totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 360
$defWs.Range("F2").Value = 113
$defWs.Range("G2").Value = 110
$defWs.Range("J2").Value = 64
$defWs.Range("O2").Value = 53
$defWs.Range("P2").Value = 40

$defWs.Range("B3").Value = 17
$defWs.Range("C3").Value = 390
$defWs.Range("D3").Value = 14
$defWs.Range("E3").Value = 44
$defWs.Range("F3").Value = 241
$defWs.Range("G3").Value = 81
$defWs.Range("H3").Value = 55
$defWs.Range("I3").Value = 131
$defWs.Range("J3").Value = 97
$defWs.Range("L3").Value = 656
$defWs.Range("M3").Value = 439
$defWs.Range("Q3").Value = 1102

# ---------------------------------------------------------------------
# ST sheet - append simulated Week 17+ special-teams logs and update
# the season totals
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 62 60"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 21 27"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 24 28 16 15"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 37"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 9"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 4 0 15 0 8"

$stWs.Range("B2").Value = 165
$stWs.Range("D2").Value = 109
$stWs.Range("F2").Value = 66
$stWs.Range("G2").Value = 63
$stWs.Range("H2").Value = 6
$stWs.Range("J2").Value = 20
$stWs.Range("K2").Value = 20
$stWs.Range("B3").Value = 83

# ---------------------------------------------------------------------
# TURNS sheet - update Road row totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 11
$turnsWs.Range("E3").Value = 21

# ---------------------------------------------------------------------
# PEN sheet - update False start (row 2) and Holding (row 3) totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 34
$penWs.Range("D2").Value = 26
$penWs.Range("B3").Value = 39
